$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns for every crypto row with the
# latest scrape snapshot. D-column writes set NumberFormat to text ("@") before
# assigning the value so number-like strings (e.g. "1.001", "240.54") are stored
# as literal text instead of being auto-coerced into a numeric value, then the
# cell style is reset to "Normal" so no stray text-format style index is left
# behind (the source cells carry no explicit style).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.809.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.900.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7665"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3052"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.36"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06849"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07988"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.904.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7360"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.170"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.837.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("E17").Value = "  -4.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.884"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007704"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.76%  "

$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.135.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.881"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.245"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1281"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.027"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.510"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.270"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.064"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05253"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.242"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7254"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.22%  "

$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01911"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.781"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.213"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4407"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.38%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8349"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.879"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.583"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.719"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.046.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.39%  "

$ws.Range("E51").Value = "  -1.30%  "

